# "update año y trimestre automáticamente tabla Resultado proyectos"
#
# The sheet "3_BD_Gerencia Nacional FINAL" recalculates the reporting
# quarter/derived columns for the "Resultado proyectos" table:
#   - Column AB ("Tiempo de gestación" helper) drops by one for every
#     data row (10 -> 9).
#   - Column Z ("Tipo Discapacidad" helper) drops by one (7 -> 6) for
#     every row currently flagged 7, reflecting the same re-numbering.
#   - Row 2's "Fecha" (AH) had not been stamped yet; it now gets the
#     same processing date already present on every other row.
#   - The active window selection moves to the last-updated row (Z27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 100

for ($r = 2; $r -le $lastRow; $r++) {
    # Column AB (28): "Tiempo de gestación" -> decrement 10 to 9
    $abCell = $ws.Cells.Item($r, 28)
    if ($abCell.Value() -eq 10) {
        $abCell.Value = 9
    }

    # Column Z (26): "Tipo Discapacidad" -> decrement 7 to 6
    $zCell = $ws.Cells.Item($r, 26)
    if ($zCell.Value() -eq 7) {
        $zCell.Value = 6
    }
}

# Row 2's "Fecha" (AH, column 34) was left blank; fill it in with the
# same date already recorded for the rest of the table (2021-10-10).
$ahCell = $ws.Cells.Item(2, 34)
if ([string]::IsNullOrEmpty($ahCell.Value())) {
    $ahCell.Value = 44479
}

# Move the active selection to reflect where the update was last applied.
$ws.Range("Z27").Select() | Out-Null
